# The "Recorded By" column (G) lists the users who recorded/edited a
# session, separated by ", ". This edit reverses the order of the
# names/emails listed in that column for every data row (rows with only a
# single recorder are left untouched since reversing a single-item list is
# a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $v = $cell.Value2

    if ($v -ne $null -and $v -ne "") {
        $parts = $v -split ", "
        $n = $parts.Count

        if ($n -gt 1) {
            $reversed = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $newVal = $reversed -join ", "
            $cell.Value = $newVal
        }
    }
}
